$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.841.57"
$ws.Range("E2").Value = "  -2.04%  "

# Row 3
$ws.Range("D3").Value = "2.904.71"
$ws.Range("E3").Value = "  -2.88%  "

# Row 4
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.02"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.01%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.79"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.61%  "

# Row 7
$ws.Range("E7").Value = "  -0.25%  "

# Row 8
$ws.Range("E8").Value = "  -0.27%  "

# Row 9
$ws.Range("D9").Value = "2.905.53"
$ws.Range("E9").Value = "  -2.49%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.75"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.02%  "

# Row 11
$ws.Range("E11").Value = "  -2.28%  "

# Row 12
$ws.Range("E12").Value = "  -3.01%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000226"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.97%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.45"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.24%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.126"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.38%  "

# Row 16
$ws.Range("D16").Value = "3.387.42"
$ws.Range("E16").Value = "  -3.04%  "

# Row 17
$ws.Range("D17").Value = "60.777.82"
$ws.Range("E17").Value = "  -2.33%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.71"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.88%  "

# Row 19
$ws.Range("D19").Value = "2.911.87"
$ws.Range("E19").Value = "  -2.88%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "429.39"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.78%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.51"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.78%  "

# Row 22
$ws.Range("E22").Value = "  -0.53%  "

# Row 23
$ws.Range("E23").Value = "  -4.61%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.06"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.20%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.88"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.45%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.21"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.69%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.98"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.26%  "

# Row 29
$ws.Range("E29").Value = "  +3.63%  "

# Row 30
$ws.Range("E30").Value = "  -0.15%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.61"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.44%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.05"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.49%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.60"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.31%  "

# Row 34
$ws.Range("E34").Value = "  -1.62%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0851"
$ws.Range("E35").Value = "  +0.46%  "

# Row 36
$ws.Range("E36").Value = "  -2.07%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.59"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.46%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.05"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.09%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.55"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.09%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.00"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.93%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.124"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.22%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.60"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.35%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.289"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.53%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.16"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -8.93%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "375.52"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.00%  "

# Row 46
$ws.Range("E46").Value = "  -1.83%  "

# Row 47
$ws.Range("D47").Value = "2.704.23"
$ws.Range("E47").Value = "  +0.81%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.11"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.81%  "

# Row 49
$ws.Range("E49").Value = "  +0.00%  "

# Row 50
$ws.Range("E50").Value = "  -7.46%  "

# Row 51
$ws.Range("E51").Value = "  -1.30%  "
